$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-9 (columns A:T), following Dr Hou's advice:
# sending/target clusters now broken out into ECs / M2 / FAPs / sCs rather than
# a single collapsed "M2" category.
$rows = @(
    @("ECs", "Ccl12", "Ccr2", "ECs", 1, 0.3333333333333333, 3.288126333333333, 9.864379, 0.05813306630866938, 0.05813306630866937, 1, 0.3333333333333333, 9.506851333333334, 28.520554, 0.04665929098818478, 0.04665929098818478, 31.25972821621844, 281.337553945966, 0.002712447656931646, 0.002712447656931645),
    @("ECs", "Ccl12", "Ccr2", "FAPs", 1, 0.3333333333333333, 3.288126333333333, 9.864379, 0.05813306630866938, 0.05813306630866937, 1, 0.3333333333333333, 0.01112833333333333, 0.033385, 0.00005461746744612846, 0.00005461746744612846, 0.03659136587944443, 0.329322292915, 0.000003175080856657377, 0.000003175080856657377),
    @("ECs", "Ccl12", "Ccr2", "M2", 1, 0.3333333333333333, 3.288126333333333, 9.864379, 0.05813306630866938, 0.05813306630866937, 3, 1, 194.1975953333333, 582.592786, 0.9531149475424379, 0.9531149475424379, 638.5462270855437, 5746.916043769894, 0.05540749444526848, 0.05540749444526847),
    @("ECs", "Ccl12", "Ccr2", "sCs", 1, 0.3333333333333333, 3.288126333333333, 9.864379, 0.05813306630866938, 0.05813306630866937, 2, 0.6666666666666666, 0.03487066666666667, 0.104612, 0.0001711440019312383, 0.0001711440019312383, 0.1146591573275555, 1.031932415948, 0.000009949125612599717, 0.000009949125612599714),
    @("M2", "Ccl12", "Ccr2", "ECs", 3, 1, 53.27393966666666, 159.821819, 0.9418669336913307, 0.9418669336913306, 1, 0.3333333333333333, 9.506851333333334, 28.520554, 0.04665929098818478, 0.04665929098818478, 506.4674243519694, 4558.206819167725, 0.04394684333125314, 0.04394684333125314),
    @("M2", "Ccl12", "Ccr2", "FAPs", 3, 1, 53.27393966666666, 159.821819, 0.9418669336913307, 0.9418669336913306, 1, 0.3333333333333333, 0.01112833333333333, 0.033385, 0.00005461746744612846, 0.00005461746744612846, 0.5928501585905553, 5.335651427314999, 0.00005144238658947109, 0.00005144238658947108),
    @("M2", "Ccl12", "Ccr2", "M2", 3, 1, 53.27393966666666, 159.821819, 0.9418669336913307, 0.9418669336913306, 3, 1, 194.1975953333333, 582.592786, 0.9531149475424379, 0.9531149475424379, 10345.67097719975, 93111.03879479773, 0.8977074530971695, 0.8977074530971694),
    @("M2", "Ccl12", "Ccr2", "sCs", 3, 1, 53.27393966666666, 159.821819, 0.9418669336913307, 0.9418669336913306, 2, 0.6666666666666666, 0.03487066666666667, 0.104612, 0.0001711440019312383, 0.0001711440019312383, 1.857697792136444, 16.719280129228, 0.0001611948763186386, 0.0001611948763186386)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowdata = $rows[$i]
    for ($j = 0; $j -lt $rowdata.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowdata[$j]
    }
}
